$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet / update the "through" date references from 07-02 to 07-03
$ws.Name = "Through 2022-07-03"
$ws.Range("I1").Value = "2022 (through 07-03)"

# Update July (row 8) 2022 value
$ws.Range("I8").Value = 23

# Update Total (row 14) 2022 value
$ws.Range("I14").Value = 829
